$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0072
$ws.Range("C2").Value = 0.9512
$ws.Range("D2").Value = 0.036
$ws.Range("E2").Value = 0.0008
$ws.Range("F2").Value = 0.0056
$ws.Range("G2").Value = 0.0056
$ws.Range("H2").Value = 0.0024
$ws.Range("I2").Value = 0.0016
$ws.Range("J2").Value = 0.0048
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0.0048
$ws.Range("M2").Value = 0.0064
$ws.Range("N2").Value = 0.0016
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0.0048
$ws.Range("Q2").Value = 0.0056
$ws.Range("R2").Value = 0.0024
$ws.Range("S2").Value = 0.0016
$ws.Range("T2").Value = 0.9448
$ws.Range("U2").Value = 0.0032
$ws.Range("V2").Value = 0.0184
$ws.Range("W2").Value = 0.0224
$ws.Range("X2").Value = 0.0048

$ws.Range("B3").Value = 0.9536
$ws.Range("C3").Value = 0.0328
$ws.Range("D3").Value = 0.0056
$ws.Range("E3").Value = 0.9568
$ws.Range("F3").Value = 0.1248
$ws.Range("G3").Value = 0.004
$ws.Range("H3").Value = 0.9616
$ws.Range("I3").Value = 0.0328
$ws.Range("J3").Value = 0.0552
$ws.Range("K3").Value = 0.0048
$ws.Range("L3").Value = 0.2056
$ws.Range("M3").Value = 0.0736
$ws.Range("N3").Value = 0.0008
$ws.Range("O3").Value = 0.9992
$ws.Range("P3").Value = 0.9216
$ws.Range("Q3").Value = 0.0024
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 0.0016
$ws.Range("T3").Value = 0.0024
$ws.Range("U3").Value = 0.9944
$ws.Range("V3").Value = 0.0072
$ws.Range("W3").Value = 0.0016
$ws.Range("X3").Value = 0.0016

$ws.Range("B4").Value = 0.0376
$ws.Range("C4").Value = 0.0056
$ws.Range("D4").Value = 0.9552
$ws.Range("E4").Value = 0.04
$ws.Range("F4").Value = 0.0024
$ws.Range("G4").Value = 0.9616
$ws.Range("H4").Value = 0.0336
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0.0008
$ws.Range("K4").Value = 0.0048
$ws.Range("L4").Value = 0.0288
$ws.Range("M4").Value = 0.0016
$ws.Range("N4").Value = 0.9952
$ws.Range("O4").Value = 0.0008
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0.9976
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0.0048
$ws.Range("U4").Value = 0.0008
$ws.Range("V4").Value = 0.0312
$ws.Range("W4").Value = 0.972
$ws.Range("X4").Value = 0.9928

$ws.Range("B5").Value = 0.0016
$ws.Range("C5").Value = 0.0104
$ws.Range("D5").Value = 0.0024
$ws.Range("E5").Value = 0.0024
$ws.Range("F5").Value = 0.8672
$ws.Range("G5").Value = 0.0288
$ws.Range("H5").Value = 0.0024
$ws.Range("I5").Value = 0.9656
$ws.Range("J5").Value = 0.9392
$ws.Range("K5").Value = 0.9904
$ws.Range("L5").Value = 0.76
$ws.Range("M5").Value = 0.9176
$ws.Range("N5").Value = 0.0024
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 0.0736
$ws.Range("Q5").Value = 0.992
$ws.Range("R5").Value = 0
$ws.Range("S5").Value = 0.9968
$ws.Range("T5").Value = 0.048
$ws.Range("U5").Value = 0.0016
$ws.Range("V5").Value = 0.9432
$ws.Range("W5").Value = 0.0032
$ws.Range("X5").Value = 0.0008
